# PictYours/Documents/Contexte.docx
# "Simple, rapide à prendre en main et élégant, l’application est intuitive..."
#                                   -> "... et élégante, l’application ..."
# (gender agreement fix: "élégant" -> "élégante")

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("élégant, l’application", $true, $false, $false, $false, $false, `
              $true, 1, $false, "élégante, l’application", 2)
